# Apply daily odds update to the "Jogos do Dia" Betfair Back/Lay workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Australian A-League Men: Melbourne City vs Macarthur FC)
$ws.Range("F2").Value = 1.69
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 5.5
$ws.Range("I2").Value = 5.8
$ws.Range("L2").Value = 1.36
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.28
$ws.Range("Q2").Value = 1.81
$ws.Range("R2").Value = 1.43
$ws.Range("T2").Value = 1.83
$ws.Range("V2").Value = 1.2
$ws.Range("W2").Value = 2.42
$ws.Range("Y2").Value = 22
$ws.Range("AI2").Value = 80
$ws.Range("AN2").Value = 9.4

# Row 4 (Portuguese Primeira Liga: Guimaraes vs Sporting Lisbon)
$ws.Range("J4").Value = 4.7
$ws.Range("AG4").Value = 42
